$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,16,17,22,23,24,26,29,32,37,38,43,44,45,47,50,53,58,59,64,65,66,68,71,74,79,80,85,86,87,88,89,90,91,99,100,105,106,107,108,109,110,111,119,120,125,126,127,128,129,130,131,139,140,145,146,147,148,149,150,151,159,160,165,166,167,168,169,170,171,179,180,185,186,187,189,192,195,200,201,206,207,208,210,213,216,221,222,227,228,229,231,234,237,242,243)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
}
